$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# Sheet "ALC": update a handful of existing numeric cells
# ---------------------------------------------------------------
$alc = $wb.Worksheets.Item("ALC")

# Row 11
$alc.Range("H11").Value = 839
$alc.Range("I11").Value = 839
$alc.Range("K11").Value = 839
$alc.Range("M11").Value = -699

# Row 42
$alc.Range("H42").Value = 2032.6666
$alc.Range("I42").Value = 98
$alc.Range("K42").Value = 294
$alc.Range("M42").Value = -64

# Row 137
$alc.Range("H137").Value = 1416.742
$alc.Range("I137").Value = 1018
$alc.Range("J137").Value = 2563.125
$alc.Range("K137").Value = 3054
$alc.Range("L137").Value = 7689.375
$alc.Range("M137").Value = -504
$alc.Range("N137").Value = -12789.375

# ---------------------------------------------------------------
# Sheet "CRP": update a handful of existing numeric cells
# ---------------------------------------------------------------
$crp = $wb.Worksheets.Item("CRP")

# Row 31
$crp.Range("H31").Value = 2258.2222
$crp.Range("I31").Value = 2016.1428
$crp.Range("J31").Value = 2656.9412
$crp.Range("K31").Value = 2016.1428
$crp.Range("L31").Value = 2656.9412
$crp.Range("M31").Value = -1721.1428
$crp.Range("N31").Value = -3246.9412

# Row 34
$crp.Range("H34").Value = 2258.2222
$crp.Range("I34").Value = 2016.1428
$crp.Range("J34").Value = 2656.9412
$crp.Range("K34").Value = 2016.1428
$crp.Range("L34").Value = 2656.9412
$crp.Range("M34").Value = -1814.1428
$crp.Range("N34").Value = -3060.9412

# ---------------------------------------------------------------
# Sheet "WVR": populate previously-empty H:N cells for rows 119-141
# (row 134 already had data and is left untouched)
# ---------------------------------------------------------------
$wvr = $wb.Worksheets.Item("WVR")

function Set-Row($Sheet, $Row, $H, $I, $J, $K, $L, $M, $N) {
    if ($null -ne $H) { $Sheet.Range("H$Row").Value = $H }
    if ($null -ne $I) { $Sheet.Range("I$Row").Value = $I }
    if ($null -ne $J) { $Sheet.Range("J$Row").Value = $J }
    if ($null -ne $K) { $Sheet.Range("K$Row").Value = $K }
    if ($null -ne $L) { $Sheet.Range("L$Row").Value = $L }
    if ($null -ne $M) { $Sheet.Range("M$Row").Value = $M }
    if ($null -ne $N) { $Sheet.Range("N$Row").Value = $N }
}

Set-Row $wvr 119 146924.5  0           146924.5  0           146924.5  $null        -156600.5
Set-Row $wvr 120 26306.334 0           26306.334 0           26306.334 $null        -35982.334
Set-Row $wvr 121 33183.547 0           33183.547 0           33183.547 $null        -36677.547
Set-Row $wvr 122 8930833   10418468    5685085   31255404    17055255  -31252954    -17060155
Set-Row $wvr 123 28860.5   0           28860.5   0           28860.5   $null        -38660.5
Set-Row $wvr 124 66000     0           66000     0           66000     $null        -75820
Set-Row $wvr 125 0         0           0         0           0         $null        $null
Set-Row $wvr 126 5269.483  5869.2      1521.25   17607.6     4563.75   -15137.6     -9503.75
Set-Row $wvr 127 17770     0           17770     0           17770     $null        -27690
Set-Row $wvr 128 0         0           0         0           0         $null        $null
Set-Row $wvr 129 0         0           0         0           0         $null        $null
Set-Row $wvr 130 0         0           0         0           0         $null        $null
Set-Row $wvr 131 22250     0           22250     0           22250     $null        -32330
Set-Row $wvr 132 1021.1579 721.6591    2034.8462 2164.9773   6104.5386 365.0227     $null
Set-Row $wvr 133 41900.5   0           41900.5   0           41900.5   $null        -52020.5
Set-Row $wvr 135 49000     0           49000     0           49000     $null        -59140
Set-Row $wvr 136 1321.6786 1321.6786   0         3965.0358   0         -1415.0358   $null
Set-Row $wvr 137 49846.668 0           49846.668 0           49846.668 $null        -60046.668
Set-Row $wvr 138 82150     0           82150     0           82150     $null        -92430
Set-Row $wvr 139 69313.336 0           69313.336 0           69313.336 $null        -79593.336
Set-Row $wvr 140 50214.5   0           50214.5   0           50214.5   $null        -60574.5
Set-Row $wvr 141 50000     0           50000     0           50000     $null        -60360

$wb.Save()
